$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(2000, "walkingToRunning", -0.443913996219635, 6.598052978515625, 6.372174263000488, 2.136175827651134, 0.4808972786212289, -4.383165293726427),
    @(2100, "walkingToRunning", 5.017239570617676, -37.66963577270508, 0.8286991119384766, 2.217999173306856, 0.9689598110900888, -3.176927756989126),
    @(2200, "walkingToRunning", -11.97968292236328, -29.96435546875, -18.5664291381836, -0.9483953673264298, 0.4404845237731658, 0.1736381053924823),
    @(2300, "walkingToRunning", 1.253250598907471, 2.138059854507446, 19.65131759643555, -6.509960667840401, -3.707480655319409, 5.622946441858697),
    @(2400, "walkingToRunning", -0.3477138876914978, -27.24154853820801, 47.39388656616211, -4.197874480965453, 6.619404304986809, 1.492906965058372),
    @(2500, "walkingToRunning", -66.68949127197266, -16.29855728149414, -26.75592613220215, 0.2072237714611235, 1.827899620450818, -3.208558189457861),
    @(2600, "walkingToRunning", -5.153336048126221, 1.702579855918884, 10.3774242401123, 4.348754659466381, -6.65085848446549, -1.474546417765225),
    @(2700, "walkingToRunning", -7.020784378051758, -33.00114822387695, 7.655538558959961, 6.596046351838384, -1.555269931239627, -1.644609145905786),
    @(2800, "walkingToRunning", 6.488995552062988, -21.21427917480469, 29.90091133117676, 1.523844735375753, 2.662038494800703, -0.68764657809816),
    @(2900, "walkingToRunning", 0.7014479041099548, 5.056992530822754, 16.41212844848633, -5.79685175829921, 0.4115632114739198, 5.635812402799219)
)

$startRow = 22
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
